# Update "PERIOD TO EXPIRE" (col H) and "LAST UPDATE" (col I) on the
# Training Dashboard sheet for rows 3-13 to reflect progress as of 04-Nov-2025.
#
# Column H values each drop by 1 day (one less day left until expiry).
# Column I ("LAST UPDATE") moves from 03-Nov-2025 to 04-Nov-2025.
#
# Writing the literal text "04-Nov-2025" straight into a General-formatted
# cell makes Excel auto-recognize it as a date and convert the cell to a
# date-serial number, which is not what we want (these are plain text
# labels in this sheet, not real dates). To keep the cell a genuine text
# value, we instead write it as a formula that evaluates to the text
# string, then immediately convert that formula to its literal value via
# Copy / PasteSpecial (values only). That collapses the cell back down to
# a plain text value without touching NumberFormat (so the cell's existing
# style/formatting is left completely untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$xlPasteValues = -4163

for ($row = 3; $row -le 13; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # Column H - PERIOD TO EXPIRE
    $hCell.Value2 = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # Column I - LAST UPDATE
    $iCell.Formula = '="04-Nov-2025"'
    $iCell.Copy()
    $iCell.PasteSpecial($xlPasteValues)
}

$excel.CutCopyMode = $false
